$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resolve duplicate variable name: "BUYED" -> "BOUGHT"
$ws.Range("A28").Value = "BOUGHT"
$ws.Range("B28").Value = "BOUGHT: "

# Changed position of exit button (reflected by moving the active selection)
$ws.Range("B28").Select()
